$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C9").Value = "Array sort and Type Conversion"
